# Change the instructions' final bullet so the exported heatmap filename
# "ex.svg" becomes "ratFlatmap.svg". The single run that currently holds
# the whole sentence is split into four runs (matching the target
# revision): "View " | "your ... as \u201c" | "ratFlatmap.svg" | "\u201d".

$d = $word.ActiveDocument

$oldName = "ex.svg"
$newName = "ratFlatmap.svg"
$splitAfter = "View "   # first run boundary, right after this literal text

# 1) Locate the filename occurrence ("ex.svg") so we don't depend on a
#    hard-coded paragraph index or character offsets.
$hit = $d.Content
$found = $hit.Find.Execute($oldName, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '$oldName' in the document"
}
$nameStart = $hit.Start
$nameEnd = $hit.End

# Paragraph that holds the filename, and its extent (excluding the
# trailing paragraph mark).
$para = $hit.Paragraphs(1)
$paraStart = $para.Range.Start
$paraEnd = $para.Range.End - 1

if (($para.Range.Text).IndexOf($splitAfter) -ne 0) {
    throw "Paragraph does not start with expected text '$splitAfter'"
}

# 2) Swap the filename text in first -- while the paragraph is still a
#    single run -- so the replacement lands cleanly in the text stream.
$nameRange = $d.Range($nameStart, $nameEnd)
$nameRange.Text = $newName
$newNameEnd = $nameStart + $newName.Length
$paraEnd = $para.Range.End - 1

# 3) Work out where the "View " / filename-quote boundaries fall now
#    that the text has its final length.
$viewSplit = $paraStart + $splitAfter.Length

# 4) Force the run to split at each boundary by toggling a character
#    property on the remainder of the paragraph (on, then back off) --
#    this leaves the text/formatting untouched but breaks the run into
#    separate <w:r> elements at each boundary.
$splitPoints = @($viewSplit, $nameStart, $newNameEnd)
foreach ($pos in $splitPoints) {
    $tail = $d.Range($pos, $paraEnd)
    $tail.Font.Bold = $true
    $tail.Font.Bold = $false
}

Write-Output "Updated paragraph: $($para.Range.Text)"
